# Add a new "2022" column (column P) to the indicator table on sheet1,
# mirroring the formatting of the existing "2021" column (O) but with the
# thousands-separated number format (#,##0.0) instead of the plain 0.0 one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numFmt = "#,##0.0"

# --- Row 4: header year ---------------------------------------------------
$ws.Range("O4").Copy()
$ws.Range("P4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("P4").Value = 2022

# --- Data rows 5-16 --------------------------------------------------------
# row -> value (numbers as numbers, the "-" placeholder as text)
$rows = @(
    @{ Row = 5;  Value = 1188.7 },
    @{ Row = 6;  Value = 263.89999999999998 },
    @{ Row = 7;  Value = 263.2 },
    @{ Row = 8;  Value = 12.4 },
    @{ Row = 9;  Value = "-" },
    @{ Row = 10; Value = 93 },
    @{ Row = 11; Value = 171.5 },
    @{ Row = 12; Value = 220.6 },
    @{ Row = 13; Value = 159.30000000000001 },
    @{ Row = 14; Value = 1.7 },
    @{ Row = 15; Value = "-" },
    @{ Row = 16; Value = 3.1 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $src = $ws.Range("O$r")
    $dst = $ws.Range("P$r")

    # Start from the same look as the corresponding "2021" (column O) cell ...
    $src.Copy()
    $dst.PasteSpecial(-4122)   # xlPasteFormats

    # ... then set the new value and the new number format ...
    $dst.Value = $item.Value
    $dst.NumberFormat = $numFmt

    # ... right-align, and drop the vertical centering the source column
    # used (the new column keeps the default/bottom vertical alignment).
    $dst.HorizontalAlignment = -4152   # xlRight
    $dst.VerticalAlignment = -4107     # xlBottom
}

$excel.CutCopyMode = $false

# --- Selection --------------------------------------------------------------
$ws.Range("Q7").Select()
